# Fixed naive component forecaster bug - Presentation state 11.02.
# Updates the y_0_forecast (column C) and y_1_forecast (column E) values
# for rows 2-39 on the active sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newC = @(
    3.123541145015474,
    3.959010658874851,
    4.722695063536686,
    6.739021039846627,
    2.619839412265601,
    -0.7919564768266385,
    1.877689851450803,
    2.705004599189187,
    1.110374544249249,
    2.267566233338814,
    2.688433258834588,
    1.014079695989589,
    3.013853578092252,
    1.331295149770684,
    0.04589006555719699,
    0.009546395482029624,
    0.8709390141433015,
    0.7652063367885598,
    2.267579219134386,
    3.146753122914103,
    1.769033835366818,
    -4.774715709990263,
    1.95493704440024,
    3.478075069442799,
    1.232342134690434,
    0.2542811494408159,
    1.519778766382096,
    1.469441753880329,
    1.638203081492495,
    2.268697431234346,
    1.984425467899631,
    0.6066448776129052,
    -4.243076347305386,
    1.438499295329754,
    1.906593537051537,
    0.08348019664223827,
    -0.214505326882275,
    0.1651547428133782
)

$newE = @(
    4.320516327661528,
    3.433494243648449,
    7.855477094481422,
    -0.7259153295281151,
    -0.17790865651377,
    2.233697987812078,
    2.959667200710037,
    0.5406927319912658,
    3.29487077883559,
    3.15890982365572,
    3.239674285955152,
    2.571626871154176,
    2.652928973511215,
    -0.3738725857433511,
    1.021287096146906,
    0.5726247744375135,
    0.6379602509701376,
    2.297544413125596,
    2.825914290412324,
    2.510325059131513,
    -1.563640406432543,
    1.685921024959058,
    3.595026567604331,
    2.123512403013161,
    1.023960954496861,
    1.414810393331356,
    0.3626330124320232,
    1.339091979913909,
    1.006355688239569,
    3.257368055312471,
    -0.3299132127116078,
    0.2691274977562275,
    23.52713615747899,
    6.182015844361843,
    1.559766133975371,
    -0.3513378667146627,
    0.256124515548195,
    -0.02867520550564606
)

$startRow = 2
for ($i = 0; $i -lt $newC.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 3).Value = $newC[$i]
    $ws.Cells.Item($row, 5).Value = $newE[$i]
}
